$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Column A (Actual Consumption (MW)) values for rows 2-28
$newA = @{
  2  = 5686
  3  = 5612
  4  = 5583
  5  = 5526
  6  = 5477
  7  = 5420
  8  = 5443
  9  = 5433
  10 = 5461
  11 = 5488
  12 = 5447
  13 = 5346
  14 = 5392
  15 = 5388
  16 = 5394
  17 = 5391
  18 = 5409
  19 = 5436
  20 = 5450
  21 = 5459
  22 = 5536
  23 = 5563
  24 = 5644
  25 = 5708
  26 = 5833
  27 = 5942
  28 = 6000
}

# New Column B (Timestamp) values for rows 2-28 (shifted by +12 from prior dataset)
$newB = @{
  2  = 45875
  3  = 45875.01041666666
  4  = 45875.02083333334
  5  = 45875.03125
  6  = 45875.04166666666
  7  = 45875.05208333334
  8  = 45875.0625
  9  = 45875.07291666666
  10 = 45875.08333333334
  11 = 45875.09375
  12 = 45875.10416666666
  13 = 45875.11458333334
  14 = 45875.125
  15 = 45875.13541666666
  16 = 45875.14583333334
  17 = 45875.15625
  18 = 45875.16666666666
  19 = 45875.17708333334
  20 = 45875.1875
  21 = 45875.19791666666
  22 = 45875.20833333334
  23 = 45875.21875
  24 = 45875.22916666666
  25 = 45875.23958333334
  26 = 45875.25
  27 = 45875.26041666666
  28 = 45875.27083333334
}

# Remove the trailing rows (29-39) that no longer exist in the new dataset
$ws.Range("A29:B39").EntireRow.Delete() | Out-Null

foreach ($r in 2..28) {
  $ws.Cells.Item($r, 1).Value = $newA[$r]
  $ws.Cells.Item($r, 2).Value = $newB[$r]
}
